$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp text
$ws.Range("C1").Value = "19.12.2019 09:51"

# Update project name
$ws.Range("A2").Value = "test"

# Update columns summary text
$ws.Range("A5").Value = "10 p. 6300 X 2"

# Resize column B width.
# Target stored XML width is 39.0625 (NPOI 1/256-char units: 39.0625*256=10000).
# This COM engine snaps stored width to (pixels+5)/6 using MDW=6, so the
# nearest reachable stored width is 39.0 (off by 0.0625) when ColumnWidth
# is set to 38.166667 (which rounds to pixel count 229 -> (229+5)/6 = 39).
$ws.Columns.Item(2).ColumnWidth = 38.166667
